# Update cryptos list (price + 1h volume change) to the refreshed snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.460.46"
$ws.Range("E2").Value = "  +5.77%  "
$ws.Range("D3").Value = "2.390.93"
$ws.Range("E3").Value = "  +3.89%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "114.79"
$ws.Range("E5").Value = "  +9.83%  "
$ws.Range("D6").Value = "318.71"
$ws.Range("E6").Value = "  +2.87%  "
$ws.Range("D7").Value = "0.635"
$ws.Range("E7").Value = "  +2.75%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.626"
$ws.Range("E9").Value = "  +3.61%  "
$ws.Range("D10").Value = "42.27"
$ws.Range("E10").Value = "  +7.32%  "
$ws.Range("D11").Value = "0.0928"
$ws.Range("E11").Value = "  +2.91%  "
$ws.Range("D12").Value = "8.65"
$ws.Range("E12").Value = "  +4.87%  "
$ws.Range("E13").Value = "  +2.76%  "
$ws.Range("E14").Value = "  +2.77%  "
$ws.Range("D15").Value = "15.92"
$ws.Range("E15").Value = "  +3.74%  "
$ws.Range("D16").Value = "2.754.96"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D17").Value = "2.389.22"
$ws.Range("E17").Value = "  +4.25%  "
$ws.Range("D18").Value = "45.510.17"
$ws.Range("E18").Value = "  +6.38%  "
$ws.Range("D19").Value = "7.48"
$ws.Range("E19").Value = "  +2.31%  "
$ws.Range("E20").Value = "  +2.96%  "
$ws.Range("D21").Value = "13.63"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("D22").Value = "74.83"
$ws.Range("E22").Value = "  +2.15%  "
$ws.Range("D23").Value = "3.56"
$ws.Range("E23").Value = "  +3.49%  "
$ws.Range("D24").Value = "264.21"
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").Value = "2.35"
$ws.Range("E25").Value = "  +6.19%  "
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").Value = "7.76"
$ws.Range("E27").Value = "  +5.71%  "
$ws.Range("D28").Value = "11.34"
$ws.Range("E28").Value = "  +4.24%  "
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("D30").Value = "39.66"
$ws.Range("E30").Value = "  +9.53%  "
$ws.Range("D31").Value = "0.0988"
$ws.Range("E31").Value = "  +15.34%  "
$ws.Range("D32").Value = "22.76"
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("D33").Value = "172.76"
$ws.Range("E33").Value = "  +4.50%  "
$ws.Range("E34").Value = "  +11.21%  "
$ws.Range("E35").Value = "  +1.56%  "
$ws.Range("D36").Value = "4.98"
$ws.Range("E36").Value = "  +9.80%  "
$ws.Range("E37").Value = "  +6.52%  "
$ws.Range("D38").Value = "4.14"
$ws.Range("E38").Value = "  +14.55%  "
$ws.Range("D39").Value = "3.08"
$ws.Range("E39").Value = "  +9.70%  "
$ws.Range("D40").Value = "0.0364"
$ws.Range("E40").Value = "  +4.83%  "
$ws.Range("D41").Value = "1.77"
$ws.Range("E41").Value = "  +12.31%  "
$ws.Range("D42").Value = "0.241"
$ws.Range("E42").Value = "  +6.66%  "
$ws.Range("D43").Value = "13.65"
$ws.Range("E43").Value = "  +11.54%  "
$ws.Range("D44").Value = "100.14"
$ws.Range("E44").Value = "  -8.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "71.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("D46").Value = "87.11"
$ws.Range("E46").Value = "  +12.83%  "
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("D48").Value = "5.83"
$ws.Range("E48").Value = "  +13.68%  "
$ws.Range("D49").Value = "115.96"
$ws.Range("E49").Value = "  +4.76%  "
$ws.Range("D50").Value = "9.48"
$ws.Range("E50").Value = "  +9.80%  "
$ws.Range("E51").Value = "  +10.67%  "
